$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("展览")
$sheet1.Range("F2").Value = 144
$sheet1.Range("F5").Value = 1240
$sheet1.Range("F6").Value = 17519
$sheet1.Range("F10").Value = 6647
$sheet1.Range("F11").Value = 675
$sheet1.Range("F13").Value = 99
$sheet1.Range("F17").Value = 162
$sheet1.Range("F24").Value = 947
$sheet1.Range("F26").Value = 5117
$sheet1.Range("F29").Value = 11795

$sheet4 = $wb.Worksheets.Item("全部类型")
$sheet4.Range("F2").Value = 144
$sheet4.Range("F5").Value = 1240
$sheet4.Range("F6").Value = 17519
$sheet4.Range("F10").Value = 6647
$sheet4.Range("F11").Value = 675
$sheet4.Range("F13").Value = 99
$sheet4.Range("F17").Value = 162
$sheet4.Range("F24").Value = 947
$sheet4.Range("F26").Value = 5117
$sheet4.Range("F30").Value = 11795
